$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.303.49'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.931.42'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.34%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9968'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.46%  '
$ws.Range("E5").Value = '  +6.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '247.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.40%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9955'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.54%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '28.40'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.98%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3226'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07107'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7903'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07997'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.34%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.930.98'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.374'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.74'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.74'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.295.29'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '253.61'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008028'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.805'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.189.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9974'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9975'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.835'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.587'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.71%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.07'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1364'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.323'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.11'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.375'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.526'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.83%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.441'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.146'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.84%  '
$ws.Range("E34").Value = '  -0.79%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.296'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7527'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.760'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01966'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.803'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '78.40'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.89%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.409'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4532'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.996'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9963'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8352'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.59'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.58%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.884'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.69%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.551'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '992.23'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +13.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.43'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1191'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.17%  '
